# Update vm_pu.xlsx results: slack-bus voltage setpoint changed from 1.05 pu to
# 1.02 pu (column B) for the 380 kV case, and the power-flow was rerun, producing
# new voltage-magnitude results in columns C..F and I..N for data rows 2-25
# (time steps 0-23, column A). Column G (the slack bus, always 1 pu) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in sheet order) whose values change for every data row.
# Column G (index 7) is intentionally omitted - it always stays at 1.
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)   # B, C, D, E, F, I, J, K, L, M, N

# One row of new values per data row (rows 2-25 <-> time steps 0-23),
# each row ordered to match $cols above: B, C, D, E, F, I, J, K, L, M, N
$data = @(
    @(1.02, 1.033354359792799, 1.049318015226544, 1.032368160686709, 1.054294154009583, 1.036834830916368, 1.03847954584039, 1.052075407958449, 1.035173593154299, 1.057037761169961, 1.039954305910752),
    @(1.02, 1.036298739072682, 1.050445439783378, 1.034948246063531, 1.056032302098939, 1.037126933458767, 1.0410565108731, 1.053014240461873, 1.037557704991747, 1.058586760831399, 1.042534930529391),
    @(1.02, 1.038191584455844, 1.051169632098825, 1.036606342287089, 1.057149188867206, 1.0373125503112, 1.042711968352587, 1.053615830106688, 1.039088774573972, 1.059580820264786, 1.044192738948382),
    @(1.02, 1.038984461140195, 1.051472828781523, 1.037300752047015, 1.057616896961679, 1.037389779236518, 1.043405121845004, 1.053867347271733, 1.039729726461478, 1.059996785779968, 1.044886876798264),
    @(1.02, 1.039117422434887, 1.051523663933549, 1.037417192980441, 1.057695320950224, 1.037402699369124, 1.043521343269362, 1.053909497131911, 1.039837188205095, 1.060066515699581, 1.045003263270376),
    @(1.02, 1.038202190099738, 1.051173688325544, 1.036615631355324, 1.057155445548895, 1.037313585398072, 1.042721241210028, 1.053619196329779, 1.039097349559703, 1.059586385978218, 1.044202024974344),
    @(1.02, 1.034352038305399, 1.049700149649811, 1.033242515891546, 1.054883210363554, 1.036934255697857, 1.039352975679201, 1.052393924644659, 1.035981761315669, 1.057562983118039, 1.040828976120118),
    @(1.02, 1.027469004386942, 1.047061804974451, 1.027208057510794, 1.050817633951864, 1.036239479162798, 1.033322252462636, 1.050188714999409, 1.030399585024821, 1.053932620460531, 1.034789688584852),
    @(1.02, 1.022808519849014, 1.045273409934123, 1.023119396528683, 1.048063285465225, 1.035758060096534, 1.029232753429166, 1.048686199644139, 1.026611656177557, 1.051466329870769, 1.030694381993682),
    @(1.02, 1.020772190595248, 1.044491696130852, 1.02133229125858, 1.046859613534259, 1.035545156165216, 1.027444474316493, 1.048027602408396, 1.02495463597778, 1.050386920111836, 1.028903563319551),
    @(1.02, 1.020012946159641, 1.044200203625969, 1.020665878058894, 1.04641080938489, 1.035465395628221, 1.026777501350523, 1.04778174090324, 1.024336526577068, 1.049984204152251, 1.028235643175466),
    @(1.02, 1.020175937947128, 1.044262781321615, 1.020808945296387, 1.046507157470878, 1.035482535412089, 1.026920694291116, 1.047834535104404, 1.024469233116919, 1.050070669221841, 1.028379039466459),
    @(1.02, 1.020709490229021, 1.044467624472033, 1.021277258993297, 1.046822550362143, 1.035538577038393, 1.027389398316503, 1.0480073046566, 1.024903596810908, 1.050353667985014, 1.028848409105324),
    @(1.02, 1.021037846921597, 1.044593684605029, 1.021565454727052, 1.047016646612574, 1.035573015900813, 1.027677818088274, 1.048113589965868, 1.02517087289947, 1.050527796155938, 1.029137238466252),
    @(1.02, 1.022943269219101, 1.045325132928852, 1.023237641127761, 1.048142932686328, 1.035772095259059, 1.029351058684662, 1.048729737729437, 1.026721265021373, 1.051537720523398, 1.030812855256211),
    @(1.02, 1.024133511562827, 1.045781969712469, 1.024282025002245, 1.048846436004158, 1.035895774743208, 1.030395885382367, 1.049114070417987, 1.027689217868115, 1.05216810918657, 1.03185916572767),
    @(1.02, 1.024826000872471, 1.04604772989483, 1.024889592873924, 1.049255717753121, 1.035967486617617, 1.031003633431789, 1.049337475365149, 1.02825219158016, 1.052534699016821, 1.032467776849032),
    @(1.02, 1.025061826789386, 1.046138228519874, 1.025096488678994, 1.049395094011868, 1.035991866184955, 1.031210577480336, 1.049413520928482, 1.028443879548857, 1.052659510674245, 1.032675014781864),
    @(1.02, 1.024005992544761, 1.045733028571823, 1.024170139043901, 1.048771066744346, 1.035882549478894, 1.030283960006543, 1.049072914979311, 1.027585533572197, 1.052100589115696, 1.031747081404974),
    @(1.02, 1.020552452262233, 1.044407334673143, 1.021139424910418, 1.04672972254958, 1.035522092993287, 1.027251452723354, 1.047956462485685, 1.024775760551016, 1.050270381304316, 1.02871026761361),
    @(1.02, 1.018364464108368, 1.04356727047796, 1.019218786766936, 1.0454363470848, 1.035291528444112, 1.025328970187953, 1.047247378688485, 1.022993948999581, 1.049109361506008, 1.026785054932497),
    @(1.02, 1.019525970179093, 1.044013234758958, 1.020238418044732, 1.046122945843247, 1.035414131309034, 1.026349647828505, 1.047623962271705, 1.023939992456208, 1.049725832663169, 1.027807182052372),
    @(1.02, 1.024063618328801, 1.045755145143714, 1.02422070044206, 1.048805126140904, 1.035888526729426, 1.030334539457011, 1.049091513727896, 1.027632389033277, 1.052131101941924, 1.031797732684063),
    @(1.02, 1.029260674468599, 1.047748974081343, 1.028779329902022, 1.051876241332486, 1.036422270580446, 1.034893138791115, 1.050764418637512, 1.03185408786116, 1.056213158537629, 1.036362805752101)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $values[$j]
    }
}

Write-Output "Updated rows 2-25 (B,C,D,E,F,I,J,K,L,M,N) for the 380 kV / 1.02 pu case."
